# Update countries & provincias Spain
# Applies the daily data refresh to the "Pais" worksheet:
#  - updates the "last updated" timestamp
#  - updates case/death numbers for several countries
#  - because some countries' total case counts now fall in a different order,
#    a handful of adjacent rows swap which country name they show (the rows'
#    numeric columns are rewritten so the visible result is the same as a
#    reordered table)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "data actualizada" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 10:42"

function Set-RowData {
    param($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes)
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Rusia - updated totals
Set-RowData 7 "Rusia" 1237504 11615 988576 227265 0 188 21663

# Turquia / Filipinas swap places (Filipinas overtakes Turquia)
Set-RowData 22 "Filipinas" 326833 2093 273313 47655 0 25 5865
Set-RowData 23 "Turquia" 326046 0 286370 31178 0 0 8498

# Indonesia - updated totals
Set-RowData 25 "Indonesia" 311176 4056 236437 63365 0 121 11374

# Singapur - updated totals
Set-RowData 60 "Singapur" 57830 11 57597 206 0 0 27

# Sudan / Eslovaquia swap places (Eslovaquia overtakes Sudan)
Set-RowData 97 "Eslovaquia" 13812 320 5027 8730 0 0 55
Set-RowData 98 "Sudan" 13653 0 6764 6053 0 0 836

# Lituania - updated totals
Set-RowData 125 "Lituania" 5366 81 2546 2721 0 5 99

# Estonia - updated totals
Set-RowData 141 "Estonia" 3659 42 2806 786 0 0 67

# Burkina Faso / Uruguay / Letonia reshuffle (Letonia overtakes both)
Set-RowData 156 "Letonia" 2194 68 1322 832 0 1 40
Set-RowData 157 "Burkina Faso" 2184 0 1420 705 0 0 59
Set-RowData 158 "Uruguay" 2155 0 1862 245 0 0 48

# Montserrat / Islas Malvinas swap places
Set-RowData 215 "Islas Malvinas" 13 0 13 0 0 0 0
Set-RowData 216 "Montserrat" 13 0 12 0 0 0 1
